# Adds a new test case (TC005) to Sheet1 that mirrors TC002 (row 3) but
# represents a "2 Drivers" scenario: a new "AdditionalDriver" column (AD) is
# introduced, and a new data row (row 6) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. New column AD: "AdditionalDriver"
# ---------------------------------------------------------------------
$ws.Range("AD1").Value = "AdditionalDriver"

# AD2 picks up the same "Nil" style used elsewhere in that row (copy
# value+format from a same-styled cell rather than just writing a value).
$ws.Range("H2").Copy($ws.Range("AD2"))

$ws.Range("AD3").Value = "Nil"
$ws.Range("AD4").Value = "Nil"
$ws.Range("AD5").Value = "Nil"

$ws.Columns.Item(30).ColumnWidth = 24.28515625

# Widen column AA (29) to fit the new, longer header text.
$ws.Columns.Item(29).ColumnWidth = 24.7109375

# ---------------------------------------------------------------------
# 2. New row 6 ("TC005"): same data as row 3 ("TC002") plus AD6 = "Yes"
# ---------------------------------------------------------------------
$ws.Range("A3:AC3").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "TC005"
$ws.Range("AD6").Value = "Yes"
$ws.Rows.Item(6).RowHeight = 48.75

# Hyperlink the new email cell, mirroring B3's mailto link. Hyperlinks.Add
# resets the cell style as a side effect, so re-apply the Hyperlink style
# afterwards (matches the style already carried over by the row copy above).
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:Email2@gmail.com")
$ws.Range("B6").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3. View / selection tidy-up to match the saved state
# ---------------------------------------------------------------------
$ws.Range("AD1").Select()
